# watchdoc v2 is added.
# Remove the "#Code ..." placeholder comment runs from the assignment
# slides, leaving the (now empty) paragraph behind. On the slide that
# has extra code already typed in (the "[]" literal), only the
# "#Code in this cell" run and its trailing line break are removed,
# while the "[]" run and its original formatting are preserved.

$p = $ppt.ActivePresentation

# Slide 3 : "Content Placeholder 2" paragraph 3 -> "#Code here"
$s = $p.Slides.Item(3)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Paragraphs(3, 1).Text = ""

# Slide 4 : "Content Placeholder 2" paragraph 2 -> "#Code Here"
$s = $p.Slides.Item(4)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Paragraphs(2, 1).Text = ""

# Slide 5 : "Content Placeholder 2" paragraph 2 -> "#Code in this cell" + <a:br/> + "[]"
# Only the "#Code in this cell" run and the line break must go; the
# "[]" run (plain Courier, no italic/color) stays.
$s = $p.Slides.Item(5)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)
$c = $tr.Characters($para.Start, $para.Length)
$c.Delete()
$tr2 = $s.Shapes.Item(1).TextFrame.TextRange
$para2 = $tr2.Paragraphs(2, 1)
$para2.Text = "[]"
$para2.Font.Name = "Courier"

# Slide 6 : "Content Placeholder 2" paragraph 3 -> "#Code in this cell"
$s = $p.Slides.Item(6)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Paragraphs(3, 1).Text = ""

# Slide 7 : "Content Placeholder 2" paragraph 2 -> "#Code in this cell"
$s = $p.Slides.Item(7)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Paragraphs(2, 1).Text = ""

# Slide 8 : "Content Placeholder 2" paragraph 3 -> "#Code in this cell"
# (paragraph 4, "Great Job!", is untouched)
$s = $p.Slides.Item(8)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Paragraphs(3, 1).Text = ""
